# Update the dSF (column F) values for several rows based on repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -1
$ws.Range("F3").Value = -5
$ws.Range("F5").Value = -8
$ws.Range("F8").Value = 0
$ws.Range("F11").Value = 4
$ws.Range("F18").Value = -2
$ws.Range("F20").Value = 4
